$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 208-224 ---
$ws = $wb.Worksheets.Item("PIR")
$pirStart = 208
$pirData = @(
    ,@('2026-01-30','17:24:38','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:39','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:40','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:40','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:41','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:42','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:43','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:47','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:52','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:24:57','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:02','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:07','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:12','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:17','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:22','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:27','17:00','Bathroom','No Motion','Inactive')
    ,@('2026-01-30','17:25:32','17:00','Bathroom','No Motion','Inactive')
)
$pirEnd = $pirStart + $pirData.Count - 1
$ws.Range("A" + $pirStart + ":A" + $pirEnd).NumberFormat = "@"
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $r = $pirStart + $i
    $ws.Cells.Item($r, 1).Value = $pirData[$i][0]
    $ws.Cells.Item($r, 2).Value = $pirData[$i][1]
    $ws.Cells.Item($r, 3).Value = $pirData[$i][2]
    $ws.Cells.Item($r, 4).Value = $pirData[$i][3]
    $ws.Cells.Item($r, 5).Value = $pirData[$i][4]
    $ws.Cells.Item($r, 6).Value = $pirData[$i][5]
}
$ws.Range("A" + $pirStart + ":A" + $pirEnd).Style = "Normal"

# --- Humidity sheet: append rows 138-148 ---
$ws = $wb.Worksheets.Item("Humidity")
$humStart = 138
$humData = @(
    ,@('2026-01-30','17:24:33','17:00','Bathroom','87.4%','Active')
    ,@('2026-01-30','17:24:38','17:00','Bathroom','87.4%','Active')
    ,@('2026-01-30','17:24:39','17:00','Bathroom','86.4%','Active')
    ,@('2026-01-30','17:24:41','17:00','Bathroom','87.4%','Active')
    ,@('2026-01-30','17:24:42','17:00','Bathroom','87.4%','Active')
    ,@('2026-01-30','17:24:48','17:00','Bathroom','87.3%','Active')
    ,@('2026-01-30','17:24:53','17:00','Bathroom','87.3%','Active')
    ,@('2026-01-30','17:24:58','17:00','Bathroom','86.4%','Active')
    ,@('2026-01-30','17:25:08','17:00','Bathroom','86.4%','Active')
    ,@('2026-01-30','17:25:13','17:00','Bathroom','87.3%','Active')
    ,@('2026-01-30','17:25:28','17:00','Bathroom','87.3%','Active')
)
$humEnd = $humStart + $humData.Count - 1
$ws.Range("A" + $humStart + ":A" + $humEnd).NumberFormat = "@"
$ws.Range("E" + $humStart + ":E" + $humEnd).NumberFormat = "@"
for ($i = 0; $i -lt $humData.Count; $i++) {
    $r = $humStart + $i
    $ws.Cells.Item($r, 1).Value = $humData[$i][0]
    $ws.Cells.Item($r, 2).Value = $humData[$i][1]
    $ws.Cells.Item($r, 3).Value = $humData[$i][2]
    $ws.Cells.Item($r, 4).Value = $humData[$i][3]
    $ws.Cells.Item($r, 5).Value = $humData[$i][4]
    $ws.Cells.Item($r, 6).Value = $humData[$i][5]
}
$ws.Range("A" + $humStart + ":A" + $humEnd).Style = "Normal"
$ws.Range("E" + $humStart + ":E" + $humEnd).Style = "Normal"

# --- mmWave sheet: append rows 48-48 ---
$ws = $wb.Worksheets.Item("mmWave")
$mmwStart = 48
$mmwData = @(
    ,@('2026-01-30','17:24:37','17:00','Living Room','FALL_DETECTED','EMERGENCY')
)
$mmwEnd = $mmwStart + $mmwData.Count - 1
$ws.Range("A" + $mmwStart + ":A" + $mmwEnd).NumberFormat = "@"
for ($i = 0; $i -lt $mmwData.Count; $i++) {
    $r = $mmwStart + $i
    $ws.Cells.Item($r, 1).Value = $mmwData[$i][0]
    $ws.Cells.Item($r, 2).Value = $mmwData[$i][1]
    $ws.Cells.Item($r, 3).Value = $mmwData[$i][2]
    $ws.Cells.Item($r, 4).Value = $mmwData[$i][3]
    $ws.Cells.Item($r, 5).Value = $mmwData[$i][4]
    $ws.Cells.Item($r, 6).Value = $mmwData[$i][5]
}
$ws.Range("A" + $mmwStart + ":A" + $mmwEnd).Style = "Normal"
